# Apply updated crypto symbol data to worksheet, preserving text cell formatting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '261.56'
Set-TextCell 'E2' '1.06%'
Set-TextCell 'G2' '18'
Set-TextCell 'E3' '1.17%'
Set-TextCell 'G3' '18'
Set-TextCell 'E4' '0.42%'
Set-TextCell 'G4' '18'
Set-TextCell 'D5' '0.06185'
Set-TextCell 'E5' '3.30%'
Set-TextCell 'G5' '18'
Set-TextCell 'D6' '6.687'
Set-TextCell 'E6' '0.45%'
Set-TextCell 'G6' '18'
Set-TextCell 'D7' '0.8511'
Set-TextCell 'E7' '-0.56%'
Set-TextCell 'G7' '18'
Set-TextCell 'D8' '0.9139'
Set-TextCell 'E8' '-1.07%'
Set-TextCell 'G8' '18'
Set-TextCell 'D9' '0.1409'
Set-TextCell 'E9' '1.42%'
Set-TextCell 'G9' '18'
Set-TextCell 'D10' '0.04646'
Set-TextCell 'E10' '-5.62%'
Set-TextCell 'G10' '18'
Set-TextCell 'D11' '0.07075'
Set-TextCell 'E11' '0.93%'
Set-TextCell 'G11' '18'
Set-TextCell 'D12' '0.03152'
Set-TextCell 'E12' '3.55%'
Set-TextCell 'G12' '18'
Set-TextCell 'D13' '0.09031'
Set-TextCell 'E13' '-1.09%'
Set-TextCell 'G13' '18'
Set-TextCell 'D14' '0.001532'
Set-TextCell 'E14' '0.52%'
Set-TextCell 'G14' '18'
Set-TextCell 'D15' '0.0006160'
Set-TextCell 'E15' '1.88%'
Set-TextCell 'G15' '18'
Set-TextCell 'D16' '0.005984'
Set-TextCell 'E16' '-1.50%'
Set-TextCell 'G16' '18'
Set-TextCell 'E17' '0.30%'
Set-TextCell 'G17' '18'
Set-TextCell 'D18' '3.169'
Set-TextCell 'E18' '0.82%'
Set-TextCell 'G18' '18'
Set-TextCell 'D19' '2.191'
Set-TextCell 'E19' '1.02%'
Set-TextCell 'G19' '18'
Set-TextCell 'E20' '-1.01%'
Set-TextCell 'G20' '18'
Set-TextCell 'E21' '0.89%'
Set-TextCell 'G21' '18'
Set-TextCell 'D22' '4.105'
Set-TextCell 'E22' '-0.71%'
Set-TextCell 'G22' '18'
Set-TextCell 'D23' '0.04217'
Set-TextCell 'E23' '-0.38%'
Set-TextCell 'G23' '18'
Set-TextCell 'D24' '0.001217'
Set-TextCell 'E24' '0.02%'
Set-TextCell 'G24' '18'
Set-TextCell 'E25' '-5.84%'
Set-TextCell 'G25' '18'
Set-TextCell 'E26' '0.08%'
Set-TextCell 'G26' '18'
Set-TextCell 'D27' '0.0001600'
Set-TextCell 'E27' '-6.50%'
Set-TextCell 'G27' '18'
Set-TextCell 'G28' '18'
Set-TextCell 'G29' '18'
Set-TextCell 'G30' '18'
Set-TextCell 'G31' '18'
Set-TextCell 'G32' '18'
Set-TextCell 'G33' '18'
Set-TextCell 'G34' '18'
Set-TextCell 'G35' '18'
Set-TextCell 'G36' '18'
Set-TextCell 'G37' '18'
Set-TextCell 'G38' '18'
Set-TextCell 'G39' '18'
Set-TextCell 'D40' '0.03901'
Set-TextCell 'E40' '1.66%'
Set-TextCell 'G40' '18'
Set-TextCell 'D41' '0.1113'
Set-TextCell 'E41' '-0.06%'
Set-TextCell 'G41' '18'
Set-TextCell 'D42' '0.004103'
Set-TextCell 'E42' '7.79%'
Set-TextCell 'G42' '18'
Set-TextCell 'E43' '-10.08%'
Set-TextCell 'G43' '18'
Set-TextCell 'D44' '0.01391'
Set-TextCell 'E44' '-7.74%'
Set-TextCell 'G44' '18'
Set-TextCell 'D45' '0.00005131'
Set-TextCell 'E45' '0.30%'
Set-TextCell 'G45' '18'
Set-TextCell 'E46' '0.08%'
Set-TextCell 'G46' '18'
Set-TextCell 'G47' '18'
Set-TextCell 'E48' '21.66%'
Set-TextCell 'G48' '18'
Set-TextCell 'D49' '0.00002100'
Set-TextCell 'E49' '0.08%'
Set-TextCell 'G49' '18'
Set-TextCell 'D50' '0.0002000'
Set-TextCell 'E50' '0.08%'
Set-TextCell 'G50' '18'
Set-TextCell 'G51' '18'
